# Add data for 2024-09-17
# Refresh the year-to-date (through 09-17) crime-count figures across the
# Citywide Totals, By Neighborhood, and individual per-neighborhood sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("E2").Value = 52
$ws.Range("J2").Value = 88
$ws.Range("D3").Value = 100
$ws.Range("G3").Value = 95
$ws.Range("H3").Value = 99
$ws.Range("J3").Value = 155
$ws.Range("G6").Value = 2
$ws.Range("B9").Value = 283
$ws.Range("C9").Value = 354
$ws.Range("H9").Value = 336
$ws.Range("I9").Value = 391
$ws.Range("J9").Value = 300
$ws.Range("K9").Value = 379
$ws.Range("B10").Value = 948
$ws.Range("C10").Value = 1159
$ws.Range("D10").Value = 1311
$ws.Range("E10").Value = 1626
$ws.Range("F10").Value = 1659
$ws.Range("H10").Value = 417
$ws.Range("I10").Value = 648
$ws.Range("J10").Value = 529
$ws.Range("K10").Value = 533
$ws.Range("B11").Value = 1337
$ws.Range("C11").Value = 1639
$ws.Range("D11").Value = 1800
$ws.Range("E11").Value = 2114
$ws.Range("F11").Value = 2222
$ws.Range("G11").Value = 1306
$ws.Range("H11").Value = 952
$ws.Range("I11").Value = 1304
$ws.Range("J11").Value = 1096
$ws.Range("K11").Value = 1207

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("C7").Value = 25
$ws.Range("E8").Value = 76
$ws.Range("H8").Value = 72
$ws.Range("I19").Value = 27
$ws.Range("J19").Value = 31
$ws.Range("E22").Value = 15
$ws.Range("B23").Value = 16
$ws.Range("E27").Value = 22
$ws.Range("B28").Value = 75
$ws.Range("D28").Value = 79
$ws.Range("I28").Value = 68
$ws.Range("J30").Value = 9
$ws.Range("F32").Value = 155
$ws.Range("K32").Value = 61
$ws.Range("B36").Value = 49
$ws.Range("G36").Value = 37
$ws.Range("F41").Value = 22
$ws.Range("B47").Value = 40
$ws.Range("C47").Value = 58
$ws.Range("D47").Value = 42
$ws.Range("K47").Value = 30
$ws.Range("I48").Value = 7
$ws.Range("H50").Value = 19
$ws.Range("C53").Value = 263
$ws.Range("D53").Value = 438
$ws.Range("E53").Value = 528
$ws.Range("F53").Value = 502
$ws.Range("H53").Value = 136
$ws.Range("I53").Value = 258
$ws.Range("K53").Value = 159
$ws.Range("K61").Value = 6
$ws.Range("B62").Value = 20
$ws.Range("H62").Value = 8
$ws.Range("D74").Value = 63
$ws.Range("J74").Value = 30
$ws.Range("C77").Value = 54
$ws.Range("G77").Value = 41
$ws.Range("D78").Value = 48
$ws.Range("D80").Value = 23
$ws.Range("C83").Value = 27
$ws.Range("B87").Value = 28
$ws.Range("E87").Value = 27
$ws.Range("H87").Value = 23
$ws.Range("J95").Value = 10
$ws.Range("F96").Value = 18
$ws.Range("E98").Value = 10
$ws.Range("B99").Value = 1337
$ws.Range("C99").Value = 1639
$ws.Range("D99").Value = 1800
$ws.Range("E99").Value = 2114
$ws.Range("F99").Value = 2222
$ws.Range("G99").Value = 1306
$ws.Range("H99").Value = 952
$ws.Range("I99").Value = 1304
$ws.Range("J99").Value = 1096
$ws.Range("K99").Value = 1207

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("E6").Value = 16
$ws.Range("E7").Value = 22

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("G5").Value = 1
$ws.Range("C9").Value = 35
$ws.Range("C10").Value = 54
$ws.Range("G10").Value = 41

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("C5").Value = 5
$ws.Range("C7").Value = 25

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J6").Value = 2
$ws.Range("J8").Value = 9

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("E2").Value = 3
$ws.Range("H6").Value = 36

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K7").Value = 21
$ws.Range("F8").Value = 102
$ws.Range("F9").Value = 155
$ws.Range("K9").Value = 61

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J2").Value = 5
$ws.Range("I8").Value = 11
$ws.Range("I9").Value = 27
$ws.Range("J9").Value = 31

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("G3").Value = 6
$ws.Range("B8").Value = 27
$ws.Range("B9").Value = 49
$ws.Range("G9").Value = 37

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("H3").Value = 11
$ws.Range("K8").Value = 39
$ws.Range("C9").Value = 221
$ws.Range("D9").Value = 379
$ws.Range("E9").Value = 467
$ws.Range("F9").Value = 447
$ws.Range("H9").Value = 71
$ws.Range("I9").Value = 152
$ws.Range("K9").Value = 82
$ws.Range("C10").Value = 263
$ws.Range("D10").Value = 438
$ws.Range("E10").Value = 528
$ws.Range("F10").Value = 502
$ws.Range("H10").Value = 136
$ws.Range("I10").Value = 258
$ws.Range("K10").Value = 159

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("H5").Value = 8
$ws.Range("H7").Value = 19

$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("D6").Value = 20
$ws.Range("D7").Value = 23

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("F6").Value = 19
$ws.Range("F7").Value = 22

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("H7").Value = 7
$ws.Range("B8").Value = 21
$ws.Range("E8").Value = 19
$ws.Range("B9").Value = 28
$ws.Range("E9").Value = 27
$ws.Range("H9").Value = 23

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("D5").Value = 44
$ws.Range("D6").Value = 48

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("D3").Value = 14
$ws.Range("I7").Value = 17
$ws.Range("B8").Value = 47
$ws.Range("B9").Value = 75
$ws.Range("D9").Value = 79
$ws.Range("I9").Value = 68

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("B7").Value = 38
$ws.Range("C7").Value = 38
$ws.Range("D7").Value = 33
$ws.Range("K7").Value = 12
$ws.Range("B8").Value = 40
$ws.Range("C8").Value = 58
$ws.Range("D8").Value = 42
$ws.Range("K8").Value = 30

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("C4").Value = 9
$ws.Range("C6").Value = 27

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J3").Value = 4
$ws.Range("D6").Value = 53
$ws.Range("D7").Value = 63
$ws.Range("J7").Value = 30

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("B7").Value = 15
$ws.Range("H7").Value = 4
$ws.Range("B8").Value = 20
$ws.Range("H8").Value = 8

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("B6").Value = 6
$ws.Range("B8").Value = 16

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("E7").Value = 13
$ws.Range("E8").Value = 15

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("J6").Value = 6
$ws.Range("J7").Value = 10

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("F5").Value = 16
$ws.Range("F6").Value = 18

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("I6").Value = 6
$ws.Range("I7").Value = 7

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("E6").Value = 9
$ws.Range("E7").Value = 10
